# Attendance was corrected: the "X" day column (7 March) is cleared out for
# every participant row (7-77). Clearing only the contents would leave the
# cell's previous "has data" border/alignment style behind, so we first
# copy the formatting of the already-blank neighboring column (Y, which is
# truly empty for these rows) onto X, then clear X's contents. This mirrors
# what Excel does when a user deletes the entries in that attendance column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Y7:Y77").Copy()
$ws.Range("X7:X77").PasteSpecial(-4122)
$ws.Range("X7:X77").ClearContents()
$excel.CutCopyMode = 0
